# Update the registration email address. This shared string is used both
# directly on RegistrationData!G1 and via formula reference on
# LoginData!A1 ( =RegistrationData!G1 ), so changing the source cell lets
# the dependent formula's cached value refresh automatically.
$wb = $excel.ActiveWorkbook

$wsRegistration = $wb.Worksheets.Item("RegistrationData")

$wsRegistration.Range("G1").Value = "ahmed.medhat1@testautomation.com"

# Switch the active sheet/selection from LoginData!B1 to
# RegistrationData!G2, matching the saved view state in the workbook.
$wsRegistration.Activate()
$wsRegistration.Range("G2").Select()
